# Fixed naive component forecaster bug - Presentation state 11.02.
#
# This script:
#   1. Clears cell C2 (a stray/erroneous value that should not be present).
#   2. Updates a set of y_1_forecast (column E) and y_0_forecast (column C)
#      values that shifted by floating point epsilon amounts once the
#      forecaster bug was fixed and the series was recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray value in C2 entirely (cell becomes empty).
$ws.Range("C2").ClearContents()

# Recomputed forecast values (tiny floating point corrections).
$ws.Range("E2").Value = 7.865470614547321
$ws.Range("E3").Value = -6.760862998203621
$ws.Range("C4").Value = 0.5799958470386724
$ws.Range("C6").Value = 0.5930547804883446
$ws.Range("E6").Value = -1.194610791899986
$ws.Range("E8").Value = 7.617133650412167
$ws.Range("C9").Value = 1.670328650030162
$ws.Range("E9").Value = 2.037906845818593
$ws.Range("C10").Value = 2.562791874943349
$ws.Range("C11").Value = 1.526411006965578
$ws.Range("E11").Value = 0.6601843988560452
$ws.Range("E12").Value = 1.55185774637272
$ws.Range("E14").Value = -5.866344937500012
$ws.Range("C15").Value = -2.616267413525608
$ws.Range("E15").Value = -4.982381489483368
$ws.Range("C17").Value = -1.298607950737285
$ws.Range("C18").Value = -0.994151974263302
$ws.Range("C19").Value = 1.069485063776932
$ws.Range("E19").Value = -2.110726282892139
